# Daily attendance processing - 2025-12-22 07:39:46
# Reverses the order of comma-separated "Recorded By" entries in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $val = $cell.Value()

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $revParts = $parts[($parts.Count - 1)..0]
            $newVal = [string]::Join(", ", $revParts)
            $cell.Value = $newVal
        }
    }
}
